# Update "想去人数" (people interested) counts on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 75
$ws1.Range("F3").Value = 3862
$ws1.Range("F4").Value = 2296
$ws1.Range("F5").Value = 453
$ws1.Range("F7").Value = 22
$ws1.Range("F8").Value = 187
$ws1.Range("F10").Value = 104
$ws1.Range("F11").Value = 1426
$ws1.Range("F12").Value = 251
$ws1.Range("F13").Value = 2497
$ws1.Range("F14").Value = 176

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 75
$ws4.Range("F3").Value = 3862
$ws4.Range("F4").Value = 2296
$ws4.Range("F5").Value = 453
$ws4.Range("F7").Value = 22
$ws4.Range("F9").Value = 187
$ws4.Range("F11").Value = 104
$ws4.Range("F14").Value = 1426
$ws4.Range("F15").Value = 251
$ws4.Range("F16").Value = 2497
$ws4.Range("F17").Value = 176
